{"js": "const paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text,style\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.style === \"Author\" && p.text.trim() === \"Edison Achalma\") {\n    target = p;\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not find the 'Edison Achalma' Author paragraph\");\n}\n\nconst newParagraph = target.insertParagraph(\n  \"Escuela Profesional de Econom\u00eda, Universidad Nacional de San Crist\u00f3bal de Huamanga\",\n  Word.InsertLocation.after\n);\nnewParagraph.style = \"Author\";\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Style.NameLocal -eq \"Author\" -and $p.Range.Text.Trim() -eq \"Edison Achalma\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -ne $null) {\n    $target.Range.InsertAfter(\"`rEscuela Profesional de Econom\u00eda, Universidad Nacional de San Crist\u00f3bal de Huamanga\")\n}\n"}
